# ---------------------------------------------------------------------------
# Commit: "Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3."
#
# What the XML diff actually contains (verified line-by-line against the
# before document extracted from before.docx):
#
#   1. <w:bookmarkStart>/<w:bookmarkEnd> w:id goes from
#      102417550852791244058099862269295220275 to
#      97159511220992923029476251964450097108
#   2. the w:rsidR GUID stamped on the REF-field runs goes from
#      44D588E79769B6E461515C556A9495D1 to 72AC777A722FB68861119C21941A525E
#   3. every toggle run-property value flips its literal spelling from
#      w:val="true"/w:val="false" to w:val="on"/w:val="off"
#      (<w:b>, <w:i>, <w:strike> only - same runs, same on/off state)
#   4. the child elements inside several <w:rPr> blocks (and one <w:shd>
#      table-cell stack) are simply re-ordered
#
# None of these four items changes a single piece of visible content: no
# run text changed, no font size/color/underline value changed, no bold/
# italic/strike *state* changed (true<->on and false<->off are the same
# boolean, just written differently), and the bookmark still wraps exactly
# "bookmark" / is still referenced by the same REF field. This is exactly
# what you'd expect from "moving from Apache POI 4.1.0 to 5.2.3": POI 5.x's
# OOXML writer spells toggle properties "on"/"off" instead of "true"/
# "false", re-orders some CTRPr children, and mints its bookmark ids/rsids
# from a different random generator - none of that is driven by, or
# reachable from, the Word object model (Bookmarks.Add never accepts a raw
# w:id, and Word has never let a script pick a run's rsid). Both id-like
# values are internal, writer-assigned implementation details, not
# document content.
#
# Concretely, in this COM host the serializer for toggle properties always
# writes its own fixed form (a bare element for "on" or w:val="0" for
# "off") no matter what a script sets Font.Bold/Italic/StrikeThrough to -
# so there is no COM call that can (re)produce the literal "true"/"false"
# or "on"/"off" strings either way. And touching *any* part of the
# document here renumbers every bookmark id down to a small sequential
# value, which would only trade one arbitrary id for a different
# arbitrary id while also disturbing unrelated paragraph-mark formatting.
#
# Since the underlying content (text, bold/italic/strike/underline state,
# font size, colors, table shading, bookmark target/name, field reference)
# is already byte-for-byte what the target expects, the correct
# application of this change through the object model is to leave the
# document's content exactly as it is - i.e. there is nothing here a Word
# script can or should mutate.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Touch nothing: the edit being replayed is a pure OOXML-writer/library
# version change (random bookmark id + rsid + boolean-literal spelling),
# with no corresponding, reachable Word-object-model operation and no
# actual content delta to apply.
